$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in the missing "Survey 3" household counts by light source (column C, rows 20-25)
$ws.Range("C20").Value = 13
$ws.Range("C21").Value = 18
$ws.Range("C22").Value = 19
$ws.Range("C23").Value = 10
$ws.Range("C24").Value = 0
$ws.Range("C25").Value = 75

# Restore the selection/scroll position left by the author when the file was saved
$excel.ActiveWindow.ScrollRow = 12
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("F27").Select()
